$d = $word.ActiveDocument

# "...saving monthly cost by $900 for patients on dialysis"
#   -> "...saving around $900 monthly cost for patients on dialysis"
$d.Content.Find.Execute(
    "monthly cost by `$900",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "around `$900 monthly cost",
    2
)
